$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2944.2222
$ws.Range("I86").Value = 2980
$ws.Range("J86").Value = 2899.5
$ws.Range("K86").Value = 2980
$ws.Range("L86").Value = 2899.5
$ws.Range("M86").Value = -1857
$ws.Range("N86").Value = -5145.5
$ws.Range("H89").Value = 2944.2222
$ws.Range("I89").Value = 2980
$ws.Range("J89").Value = 2899.5
$ws.Range("K89").Value = 14900
$ws.Range("L89").Value = 14497.5
$ws.Range("M89").Value = -9284
$ws.Range("N89").Value = -25729.5
$ws.Range("H101").Value = 1510.6875
$ws.Range("I101").Value = 1279.5555
$ws.Range("K101").Value = 3838.6665
$ws.Range("M101").Value = -2216.6665
$ws.Range("H125").Value = 2406
$ws.Range("I125").Value = 2500.3333
$ws.Range("K125").Value = 22502.9997
$ws.Range("M125").Value = -20042.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3892.5
$ws.Range("I61").Value = 2926.0667
$ws.Range("J61").Value = 5503.222
$ws.Range("K61").Value = 2926.0667
$ws.Range("L61").Value = 5503.222
$ws.Range("M61").Value = -2714.0667
$ws.Range("N61").Value = -5927.222
$ws.Range("H74").Value = 245746.44
$ws.Range("I74").Value = 500477.84
$ws.Range("J74").Value = 3145.0952
$ws.Range("K74").Value = 500477.84
$ws.Range("L74").Value = 3145.0952
$ws.Range("M74").Value = -499603.84
$ws.Range("N74").Value = -4893.0952
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 245746.44
$ws.Range("I77").Value = 500477.84
$ws.Range("J77").Value = 3145.0952
$ws.Range("K77").Value = 2502389.2
$ws.Range("L77").Value = 15725.476
$ws.Range("M77").Value = -2498021.2
$ws.Range("N77").Value = -24461.476
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H122").Value = 1280.9688
$ws.Range("I122").Value = 975.3214
$ws.Range("J122").Value = 3420.5
$ws.Range("K122").Value = 2925.9642
$ws.Range("L122").Value = 10261.5
$ws.Range("M122").Value = -475.9642000000003
$ws.Range("N122").Value = -15161.5
$ws.Range("H132").Value = 3230.0334
$ws.Range("I132").Value = 2248.32
$ws.Range("K132").Value = 6744.960000000001
$ws.Range("M132").Value = -4214.960000000001
$ws.Range("H136").Value = 3892.5
$ws.Range("I136").Value = 2926.0667
$ws.Range("J136").Value = 5503.222
$ws.Range("K136").Value = 8778.2001
$ws.Range("L136").Value = 16509.666
$ws.Range("M136").Value = -6228.2001
$ws.Range("N136").Value = -21609.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 60000
$ws.Range("J45").Value = 60000
$ws.Range("L45").Value = 60000
$ws.Range("N45").Value = -61616
$ws.Range("H86").Value = 108990.21
$ws.Range("I86").Value = 3929.9412
$ws.Range("J86").Value = 1002002.5
$ws.Range("K86").Value = 3929.9412
$ws.Range("L86").Value = 1002002.5
$ws.Range("M86").Value = -2806.9412
$ws.Range("N86").Value = -1004248.5
$ws.Range("H89").Value = 108990.21
$ws.Range("I89").Value = 3929.9412
$ws.Range("J89").Value = 1002002.5
$ws.Range("K89").Value = 19649.706
$ws.Range("L89").Value = 5010012.5
$ws.Range("M89").Value = -14033.706
$ws.Range("N89").Value = -5021244.5
$ws.Range("H99").Value = 2097.25
$ws.Range("I99").Value = 1611.1538
$ws.Range("K99").Value = 1611.1538
$ws.Range("M99").Value = -113.1538
$ws.Range("H134").Value = 2724.8
$ws.Range("I134").Value = 1265.8
$ws.Range("J134").Value = 10019.8
$ws.Range("K134").Value = 3797.4
$ws.Range("L134").Value = 30059.4
$ws.Range("M134").Value = -1262.4
$ws.Range("N134").Value = -35129.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34486484
$ws.Range("J31").Value = 5426.533
$ws.Range("L31").Value = 5426.533
$ws.Range("N31").Value = -6016.533
$ws.Range("H34").Value = 34486484
$ws.Range("J34").Value = 5426.533
$ws.Range("L34").Value = 5426.533
$ws.Range("N34").Value = -5830.533
$ws.Range("H86").Value = 6090.3335
$ws.Range("I86").Value = 4257.857
$ws.Range("K86").Value = 4257.857
$ws.Range("M86").Value = -3134.857
$ws.Range("H89").Value = 6090.3335
$ws.Range("I89").Value = 4257.857
$ws.Range("K89").Value = 21289.285
$ws.Range("M89").Value = -15673.285
$ws.Range("H95").Value = 12470.777
$ws.Range("J95").Value = 12470.777
$ws.Range("L95").Value = 12470.777
$ws.Range("N95").Value = -17962.777
$ws.Range("H122").Value = 112215.11
$ws.Range("I122").Value = 143991
$ws.Range("K122").Value = 431973
$ws.Range("M122").Value = -429523
$ws.Range("H131").Value = 49999
$ws.Range("J131").Value = 49999
$ws.Range("L131").Value = 49999
$ws.Range("N131").Value = -60079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 242.45454
$ws.Range("I38").Value = 216.7
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 650.0999999999999
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = -303.0999999999999
$ws.Range("N38").Value = -2194

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 30000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 30000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 30000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -30518
$ws.Range("H80").Value = 2982.6667
$ws.Range("I80").Value = 1999
$ws.Range("J80").Value = 3966.3333
$ws.Range("K80").Value = 1999
$ws.Range("L80").Value = 3966.3333
$ws.Range("M80").Value = -1001
$ws.Range("N80").Value = -5962.3333
$ws.Range("H83").Value = 2982.6667
$ws.Range("I83").Value = 1999
$ws.Range("J83").Value = 3966.3333
$ws.Range("K83").Value = 9995
$ws.Range("L83").Value = 19831.6665
$ws.Range("M83").Value = -5003
$ws.Range("N83").Value = -29815.6665
$ws.Range("H92").Value = 28500.125
$ws.Range("J92").Value = 28285.857
$ws.Range("L92").Value = 28285.857
$ws.Range("N92").Value = -32029.857
$ws.Range("H126").Value = 2528.7
$ws.Range("I126").Value = 2257.375
$ws.Range("K126").Value = 6772.125
$ws.Range("M126").Value = -4302.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2917.4285
$ws.Range("I22").Value = 1811.125
$ws.Range("J22").Value = 4392.5
$ws.Range("K22").Value = 1811.125
$ws.Range("L22").Value = 4392.5
$ws.Range("M22").Value = -1516.125
$ws.Range("N22").Value = -4982.5
$ws.Range("H27").Value = 2917.4285
$ws.Range("I27").Value = 1811.125
$ws.Range("J27").Value = 4392.5
$ws.Range("K27").Value = 1811.125
$ws.Range("L27").Value = 4392.5
$ws.Range("M27").Value = -1704.125
$ws.Range("N27").Value = -4606.5
$ws.Range("H30").Value = 15000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 15000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 15000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -15216
$ws.Range("H61").Value = 2394
$ws.Range("I61").Value = 868.4
$ws.Range("K61").Value = 868.4
$ws.Range("M61").Value = -666.4
$ws.Range("H93").Value = 5097.6665
$ws.Range("I93").Value = 4868.143
$ws.Range("J93").Value = 5243.727
$ws.Range("K93").Value = 4868.143
$ws.Range("L93").Value = 5243.727
$ws.Range("M93").Value = -3620.143
$ws.Range("N93").Value = -7739.727
$ws.Range("H113").Value = 2394
$ws.Range("I113").Value = 868.4
$ws.Range("K113").Value = 868.4
$ws.Range("M113").Value = 1301.6
$ws.Range("H132").Value = 7485.7856
$ws.Range("I132").Value = 5913.6
$ws.Range("J132").Value = 11416.25
$ws.Range("K132").Value = 17740.8
$ws.Range("L132").Value = 34248.75
$ws.Range("M132").Value = -15210.8
$ws.Range("N132").Value = -39308.75
$ws.Range("H136").Value = 4625.3237
$ws.Range("I136").Value = 3633.3076
$ws.Range("J136").Value = 5239.4287
$ws.Range("K136").Value = 10899.9228
$ws.Range("L136").Value = 15718.2861
$ws.Range("M136").Value = -8349.9228
$ws.Range("N136").Value = -20818.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 31881.5
$ws.Range("J34").Value = 31176
$ws.Range("L34").Value = 31176
$ws.Range("N34").Value = -31582
$ws.Range("H41").Value = 29999.5
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 29999.5
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 29999.5
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -30779.5
$ws.Range("H57").Value = 79999
$ws.Range("J57").Value = 79999
$ws.Range("L57").Value = 79999
$ws.Range("N57").Value = -81507
$ws.Range("H69").Value = 22000
$ws.Range("J69").Value = 22000
$ws.Range("L69").Value = 22000
$ws.Range("N69").Value = -23498
$ws.Range("H72").Value = 22000
$ws.Range("J72").Value = 22000
$ws.Range("L72").Value = 66000
$ws.Range("N72").Value = -73488
$ws.Range("H126").Value = 3758.0833
$ws.Range("I126").Value = 2142.158
$ws.Range("J126").Value = 9898.6
$ws.Range("K126").Value = 6426.474
$ws.Range("L126").Value = 29695.8
$ws.Range("M126").Value = -3956.474
$ws.Range("N126").Value = -34635.8
$ws.Range("H136").Value = 5906.8667
$ws.Range("I136").Value = 2872
$ws.Range("K136").Value = 8616
$ws.Range("M136").Value = -6066
